$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Severity value for Bug 2 (row 7) from 1 to 3
$ws.Range("B7").Value = 3

# Update the active selection to B29 (also resets the scrolled view)
$ws.Range("B29").Select() | Out-Null
